$d = $word.ActiveDocument

# 1. Update the "npm install @angular/core@15.2.1" line to append " -legacy-peer-deps"
#    (en dash, consistent with the other "-legacy-peer-deps" / "-force" lines in the doc)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "npm install @angular/core@15.2.1") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Text = "npm install @angular/core@15.2.1 " + [char]8211 + "legacy-peer-deps"

    # 2. Insert a new paragraph right after it with the "npm install -g @angular/cli" text
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Range.Text = "npm install -g @angular/cli"
}

# 3. Suppress automatic hyphenation on the "Normal" style (adds <w:suppressAutoHyphens/> to its pPr)
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.Hyphenation = $false
